# Generate Report for Handback
#
# Two handback files are renamed/regenerated:
#   668e1d8b-ca32-4470-bc2d-b3a2537e67a4.md  -> a4de48c1-632b-4595-8be7-a0ccb2a4b613.md
#   99b90781-9224-4582-ba7b-4fe81cf19a3c.md  -> ffff07ec4d44-2df2-4a84-94f6-0122f514a74a.md
#
# and both now share a single regenerated xliff pair (one for zh-cn, one for
# de-de) with refreshed timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "a4de48c1-632b-4595-8be7-a0ccb2a4b613.md"
$ov.Range("B2").Value = "e2e\a4de48c1-632b-4595-8be7-a0ccb2a4b613.md"
$ov.Range("G2").Value = "2016-08-31 11:08:59"

$ov.Range("A3").Value = "ffff07ec4d44-2df2-4a84-94f6-0122f514a74a.md"
$ov.Range("B3").Value = "e2e\ffff07ec4d44-2df2-4a84-94f6-0122f514a74a.md"
$ov.Range("G3").Value = "2016-08-31 11:08:59"

# Hyperlinks (B2/B3) need their display text refreshed. This shim's
# Hyperlink objects can't be mutated/removed individually, so rebuild the
# sheet's hyperlink set from scratch, preserving the original target URLs.
$ovLinks = $ov.Hyperlinks
$ovLinks.Delete()
$ovLinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/259ec66d00347768ed2d7338bf3de2bea2b732b2/e2e/668e1d8b-ca32-4470-bc2d-b3a2537e67a4.md", "", "", "e2e\a4de48c1-632b-4595-8be7-a0ccb2a4b613.md")
$ovLinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/259ec66d00347768ed2d7338bf3de2bea2b732b2/e2e/99b90781-9224-4582-ba7b-4fe81cf19a3c.md", "", "", "e2e\ffff07ec4d44-2df2-4a84-94f6-0122f514a74a.md")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "a4de48c1-632b-4595-8be7-a0ccb2a4b613.md"
$zh.Range("G2").Value = "a4de48c1-632b-4595-8be7-a0ccb2a4b613.09be350934fcfc3daa36d49e094959398804949d.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-31 11:08:54"
$zh.Range("I2").Value = "a4de48c1-632b-4595-8be7-a0ccb2a4b613.md"
$zh.Range("J2").Value = "a4de48c1-632b-4595-8be7-a0ccb2a4b613.09be350934fcfc3daa36d49e094959398804949d.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-31 11:09:17"

$zh.Range("A3").Value = "ffff07ec4d44-2df2-4a84-94f6-0122f514a74a.md"
$zh.Range("G3").Value = "a4de48c1-632b-4595-8be7-a0ccb2a4b613.09be350934fcfc3daa36d49e094959398804949d.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-31 11:08:54"
$zh.Range("I3").Value = "ffff07ec4d44-2df2-4a84-94f6-0122f514a74a.md"
$zh.Range("J3").Value = "a4de48c1-632b-4595-8be7-a0ccb2a4b613.09be350934fcfc3daa36d49e094959398804949d.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-31 11:09:17"

$zhLinks = $zh.Hyperlinks
$zhLinks.Delete()
$zhLinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/259ec66d00347768ed2d7338bf3de2bea2b732b2/e2e/668e1d8b-ca32-4470-bc2d-b3a2537e67a4.md", "", "", "a4de48c1-632b-4595-8be7-a0ccb2a4b613.md")
$zhLinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0783ec277ab8cb05394e204b93f07c12c3388899/e2e/668e1d8b-ca32-4470-bc2d-b3a2537e67a4.md", "", "", "a4de48c1-632b-4595-8be7-a0ccb2a4b613.md")
$zhLinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/259ec66d00347768ed2d7338bf3de2bea2b732b2/e2e/99b90781-9224-4582-ba7b-4fe81cf19a3c.md", "", "", "ffff07ec4d44-2df2-4a84-94f6-0122f514a74a.md")
$zhLinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0783ec277ab8cb05394e204b93f07c12c3388899/e2e/99b90781-9224-4582-ba7b-4fe81cf19a3c.md", "", "", "ffff07ec4d44-2df2-4a84-94f6-0122f514a74a.md")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "a4de48c1-632b-4595-8be7-a0ccb2a4b613.md"
$de.Range("G2").Value = "a4de48c1-632b-4595-8be7-a0ccb2a4b613.09be350934fcfc3daa36d49e094959398804949d.de-de.xlf"
$de.Range("H2").Value = "2016-08-31 11:08:59"
$de.Range("I2").Value = "a4de48c1-632b-4595-8be7-a0ccb2a4b613.md"
$de.Range("J2").Value = "a4de48c1-632b-4595-8be7-a0ccb2a4b613.09be350934fcfc3daa36d49e094959398804949d.de-de.xlf"
$de.Range("K2").Value = "2016-08-31 11:09:24"

$de.Range("A3").Value = "ffff07ec4d44-2df2-4a84-94f6-0122f514a74a.md"
$de.Range("G3").Value = "a4de48c1-632b-4595-8be7-a0ccb2a4b613.09be350934fcfc3daa36d49e094959398804949d.de-de.xlf"
$de.Range("H3").Value = "2016-08-31 11:08:59"
$de.Range("I3").Value = "ffff07ec4d44-2df2-4a84-94f6-0122f514a74a.md"
$de.Range("J3").Value = "a4de48c1-632b-4595-8be7-a0ccb2a4b613.09be350934fcfc3daa36d49e094959398804949d.de-de.xlf"
$de.Range("K3").Value = "2016-08-31 11:09:24"

$deLinks = $de.Hyperlinks
$deLinks.Delete()
$deLinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/259ec66d00347768ed2d7338bf3de2bea2b732b2/e2e/668e1d8b-ca32-4470-bc2d-b3a2537e67a4.md", "", "", "a4de48c1-632b-4595-8be7-a0ccb2a4b613.md")
$deLinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/db97ed87e1dfca09d4e250b1588f704b1c8d2645/e2e/668e1d8b-ca32-4470-bc2d-b3a2537e67a4.md", "", "", "a4de48c1-632b-4595-8be7-a0ccb2a4b613.md")
$deLinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/259ec66d00347768ed2d7338bf3de2bea2b732b2/e2e/99b90781-9224-4582-ba7b-4fe81cf19a3c.md", "", "", "ffff07ec4d44-2df2-4a84-94f6-0122f514a74a.md")
$deLinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/db97ed87e1dfca09d4e250b1588f704b1c8d2645/e2e/99b90781-9224-4582-ba7b-4fe81cf19a3c.md", "", "", "ffff07ec4d44-2df2-4a84-94f6-0122f514a74a.md")

Write-Host "Done."
